# Generate Report for Handoff
# Adds two new handoff entries (79f67d66-... and e54e3010-...) as new rows
# to the "Overview", "zh-cn" and "de-de" tables/sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet (File Name / Path And Name / Extension / Publish URL /
# zh-cn / de-de / Latest HO Xliff Generate Date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

# Row 4 : 79f67d66-c291-4c74-9f5e-95a13b554a22.md
$wsOverview.Range("A4").Value = "79f67d66-c291-4c74-9f5e-95a13b554a22.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/79f67d66-c291-4c74-9f5e-95a13b554a22.md", "", "", "e2e\79f67d66-c291-4c74-9f5e-95a13b554a22.md") | Out-Null
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-09-05 08:55:34"

# Row 5 : e54e3010-17cf-4717-a343-dedf557bcfd5.md
$wsOverview.Range("A5").Value = "e54e3010-17cf-4717-a343-dedf557bcfd5.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/e54e3010-17cf-4717-a343-dedf557bcfd5.md", "", "", "e2e\e54e3010-17cf-4717-a343-dedf557bcfd5.md") | Out-Null
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-09-05 08:55:34"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null
$loZh.ListRows.Add() | Out-Null

# Row 4 : 79f67d66-c291-4c74-9f5e-95a13b554a22.md
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/79f67d66-c291-4c74-9f5e-95a13b554a22.md", "", "", "79f67d66-c291-4c74-9f5e-95a13b554a22.md") | Out-Null
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "False"
$wsZh.Range("G4").Value = "79f67d66-c291-4c74-9f5e-95a13b554a22.d1e71096976d94bdce89d905b85b087c97347128.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-09-05 08:55:29"
$wsZh.Range("I4").Value = ""
$wsZh.Range("J4").Value = ""
$wsZh.Range("K4").Value = "0001-01-01 00:00:00"
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "False"
$wsZh.Range("P4").Value = ""

# Row 5 : e54e3010-17cf-4717-a343-dedf557bcfd5.md
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/e54e3010-17cf-4717-a343-dedf557bcfd5.md", "", "", "e54e3010-17cf-4717-a343-dedf557bcfd5.md") | Out-Null
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "False"
$wsZh.Range("G5").Value = "e54e3010-17cf-4717-a343-dedf557bcfd5.b6b83d25c79b4287f915840084b611e1352cf812.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-09-05 08:55:29"
$wsZh.Range("I5").Value = ""
$wsZh.Range("J5").Value = ""
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("L5").Value = ""
$wsZh.Range("M5").Value = "True"
$wsZh.Range("N5").Value = ""
$wsZh.Range("O5").Value = "False"
$wsZh.Range("P5").Value = ""

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null
$loDe.ListRows.Add() | Out-Null

# Row 4 : 79f67d66-c291-4c74-9f5e-95a13b554a22.md
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/79f67d66-c291-4c74-9f5e-95a13b554a22.md", "", "", "79f67d66-c291-4c74-9f5e-95a13b554a22.md") | Out-Null
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "False"
$wsDe.Range("G4").Value = "79f67d66-c291-4c74-9f5e-95a13b554a22.d1e71096976d94bdce89d905b85b087c97347128.de-de.xlf"
$wsDe.Range("H4").Value = "2016-09-05 08:55:34"
$wsDe.Range("I4").Value = ""
$wsDe.Range("J4").Value = ""
$wsDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "False"
$wsDe.Range("P4").Value = ""

# Row 5 : e54e3010-17cf-4717-a343-dedf557bcfd5.md
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/e54e3010-17cf-4717-a343-dedf557bcfd5.md", "", "", "e54e3010-17cf-4717-a343-dedf557bcfd5.md") | Out-Null
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "False"
$wsDe.Range("G5").Value = "e54e3010-17cf-4717-a343-dedf557bcfd5.b6b83d25c79b4287f915840084b611e1352cf812.de-de.xlf"
$wsDe.Range("H5").Value = "2016-09-05 08:55:34"
$wsDe.Range("I5").Value = ""
$wsDe.Range("J5").Value = ""
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("L5").Value = ""
$wsDe.Range("M5").Value = "True"
$wsDe.Range("N5").Value = ""
$wsDe.Range("O5").Value = "False"
$wsDe.Range("P5").Value = ""

Write-Host "Handback rows added for 79f67d66-... and e54e3010-... across Overview, zh-cn, de-de."
